# Apply updated cryptocurrency price/volume data per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.87"
$ws.Range("E2").Value = "'0.78%"
$ws.Range("D3").Value = "'29.40"
$ws.Range("E3").Value = "'7.11%"
$ws.Range("D4").Value = "'5.186"
$ws.Range("E4").Value = "'1.37%"
$ws.Range("D5").Value = "'0.05737"
$ws.Range("E5").Value = "'0.97%"
$ws.Range("D6").Value = "'6.560"
$ws.Range("E6").Value = "'0.66%"
$ws.Range("D7").Value = "'0.8591"
$ws.Range("E7").Value = "'4.84%"
$ws.Range("D8").Value = "'0.8666"
$ws.Range("E8").Value = "'1.75%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1367"
$ws.Range("E9").Value = "'2.51%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07085"
$ws.Range("E10").Value = "'1.94%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.03013"
$ws.Range("E11").Value = "'4.65%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09385"
$ws.Range("E12").Value = "'-0.12%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001535"
$ws.Range("E13").Value = "'0.61%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006002"
$ws.Range("E14").Value = "'0.37%"
$ws.Range("D15").Value = "'0.005966"
$ws.Range("E15").Value = "'-4.00%"
$ws.Range("E16").Value = "'5,225.17%"
$ws.Range("D17").Value = "'3.490"
$ws.Range("E17").Value = "'-0.58%"
$ws.Range("D18").Value = "'3.106"
$ws.Range("E18").Value = "'3.21%"
$ws.Range("D19").Value = "'2.188"
$ws.Range("E19").Value = "'-5.61%"
$ws.Range("D20").Value = "'0.3200"
$ws.Range("E20").Value = "'0.26%"
$ws.Range("D21").Value = "'0.03309"
$ws.Range("E21").Value = "'2.64%"
$ws.Range("D22").Value = "'0.1294"
$ws.Range("E22").Value = "'1.53%"
$ws.Range("D23").Value = "'3.495"
$ws.Range("E23").Value = "'-1.76%"
$ws.Range("D24").Value = "'0.04144"
$ws.Range("E24").Value = "'2.96%"
$ws.Range("D25").Value = "'0.1381"
$ws.Range("E25").Value = "'0.48%"
$ws.Range("D26").Value = "'0.001227"
$ws.Range("E26").Value = "'0.96%"
$ws.Range("E27").Value = "'11.53%"
$ws.Range("E28").Value = "'2.61%"
$ws.Range("E40").Value = "'0.89%"
$ws.Range("D41").Value = "'0.005764"
$ws.Range("E41").Value = "'-1.30%"
$ws.Range("D42").Value = "'0.1074"
$ws.Range("E42").Value = "'1.48%"
$ws.Range("D43").Value = "'0.002442"
$ws.Range("E43").Value = "'6.16%"
$ws.Range("D44").Value = "'0.009479"
$ws.Range("E44").Value = "'-2.37%"
$ws.Range("D45").Value = "'0.00005292"
$ws.Range("E45").Value = "'3.78%"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("D47").Value = "'0.05702"
$ws.Range("E47").Value = "'-43.54%"
$ws.Range("D48").Value = "'0.002282"
$ws.Range("E48").Value = "'-9.41%"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("E50").Value = "'0.04%"
